# Apply cryptos list update (Mon Aug 26 23:37:45 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Sheet, $Addr, $Text) {
    $rng = $Sheet.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "63.049.82"
Set-TextCell $ws "E2" "  -2.50%  "

Set-TextCell $ws "D3" "2.684.41"
Set-TextCell $ws "E3" "  -2.79%  "

Set-TextCell $ws "E4" "  -0.02%  "

Set-TextCell $ws "D5" "552.06"
Set-TextCell $ws "E5" "  -4.37%  "

Set-TextCell $ws "D6" "157.89"
Set-TextCell $ws "E6" "  -1.51%  "

Set-TextCell $ws "E7" "  +0.12%  "

Set-TextCell $ws "D8" "0.588"
Set-TextCell $ws "E8" "  -2.47%  "

Set-TextCell $ws "E9" "  -4.66%  "

Set-TextCell $ws "E10" "  -2.96%  "

Set-TextCell $ws "E11" "  -4.78%  "

Set-TextCell $ws "D12" "5.10"
Set-TextCell $ws "E12" "  -12.61%  "

Set-TextCell $ws "D13" "3.160.48"
Set-TextCell $ws "E13" "  -2.80%  "

Set-TextCell $ws "D14" "26.11"
Set-TextCell $ws "E14" "  -4.55%  "

Set-TextCell $ws "D15" "62.924.73"
Set-TextCell $ws "E15" "  -2.09%  "

Set-TextCell $ws "E16" "  -3.82%  "

Set-TextCell $ws "D17" "2.686.59"
Set-TextCell $ws "E17" "  -3.01%  "

Set-TextCell $ws "D18" "11.91"
Set-TextCell $ws "E18" "  -2.23%  "

Set-TextCell $ws "D19" "4.58"
Set-TextCell $ws "E19" "  -5.68%  "

Set-TextCell $ws "D20" "343.52"
Set-TextCell $ws "E20" "  -4.20%  "

Set-TextCell $ws "E21" "  -5.41%  "

Set-TextCell $ws "D22" "1.00"
Set-TextCell $ws "E22" "  +0.06%  "

Set-TextCell $ws "D23" "0.505"
Set-TextCell $ws "E23" "  -4.67%  "

Set-TextCell $ws "D24" "63.53"
Set-TextCell $ws "E24" "  -2.52%  "

Set-TextCell $ws "E25" "  -1.96%  "

Set-TextCell $ws "E27" "  -5.66%  "

Set-TextCell $ws "D28" "0.0₃0857"
Set-TextCell $ws "E28" "  -7.54%  "

Set-TextCell $ws "E29" "  -2.01%  "

Set-TextCell $ws "E30" "  -3.47%  "

Set-TextCell $ws "D31" "7.04"
Set-TextCell $ws "E31" "  -4.74%  "

Set-TextCell $ws "D32" "167.45"
Set-TextCell $ws "E32" "  -0.14%  "

Set-TextCell $ws "E33" "  +0.05%  "

Set-TextCell $ws "D34" "4.82"
Set-TextCell $ws "E34" "  -4.02%  "

Set-TextCell $ws "D35" "19.53"

Set-TextCell $ws "E36" "  -5.83%  "

Set-TextCell $ws "D37" "1.77"
Set-TextCell $ws "E37" "  -4.55%  "

Set-TextCell $ws "D38" "339.71"
Set-TextCell $ws "E38" "  -3.65%  "

Set-TextCell $ws "D39" "6.20"
Set-TextCell $ws "E39" "  -3.64%  "

Set-TextCell $ws "D40" "0.932"
Set-TextCell $ws "E40" "  -7.59%  "

Set-TextCell $ws "B41" "OKB"
Set-TextCell $ws "C41" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D41" "38.27"
Set-TextCell $ws "E41" "  -2.28%  "

Set-TextCell $ws "B42" "Filecoin"
Set-TextCell $ws "C42" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws "D42" "3.94"
Set-TextCell $ws "E42" "  -6.02%  "

Set-TextCell $ws "D43" "20.32"
Set-TextCell $ws "E43" "  -6.10%  "

Set-TextCell $ws "D44" "20.79"
Set-TextCell $ws "E44" "  -8.20%  "

Set-TextCell $ws "D45" "0.617"
Set-TextCell $ws "E45" "  -2.48%  "

Set-TextCell $ws "E47" "  -0.06%  "

Set-TextCell $ws "E48" "  +0.06%  "

Set-TextCell $ws "D49" "0.0972"
Set-TextCell $ws "E49" "  -4.20%  "

Set-TextCell $ws "D50" "129.14"
Set-TextCell $ws "E50" "  -5.48%  "

Set-TextCell $ws "D51" "2.092.26"
Set-TextCell $ws "E51" "  -2.54%  "
